$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 07:22"

# Australia (row 24): update Casos activos (D) and Recuperados (E)
$ws.Range("D24").Value = 1080
$ws.Range("E24").Value = 4770

# Japon (row 33): update Muertes hoy (F)
$ws.Range("F33").Value = 79

# Pakistan moves above Malasia (row 34 becomes Pakistan with updated stats,
# row 35 becomes Malasia with its previous stats)
$ws.Range("A34").Value = "Pakistan"
$ws.Range("B34").Value = 3861
$ws.Range("C34").Value = 95
$ws.Range("D34").Value = 259
$ws.Range("E34").Value = 3548
$ws.Range("F34").Value = 17
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 54

$ws.Range("A35").Value = "Malasia"
$ws.Range("B35").Value = 3793
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 1241
$ws.Range("E35").Value = 2490
$ws.Range("F35").Value = 102
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 62

# Nueva Zelanda (row 60): update Muertes hoy (F)
$ws.Range("F60").Value = 14

# Cuba (row 93): update Casos totales (B), Nuevos casos (C), Recuperados (E)
$ws.Range("B93").Value = 363
$ws.Range("C93").Value = 13
$ws.Range("E93").Value = 336

# Sri Lanka (row 113): update Recuperados (E), Casos criticos (G), Muertes (H)
$ws.Range("E113").Value = 134
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 6
